# Fruta / hortaliza, semanal
# Insert two new weekly records (row 168 and 169) into the Limón dataset,
# pushing all subsequent rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 168 (shifts existing 168.. down to 170..)
$ws.Range("A168:A169").EntireRow.Insert()

# New row 168: 1a amarillo
$ws.Cells.Item(168, 1).Value = 4
$ws.Cells.Item(168, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(168, 3).Value = "Los Lagos"
$ws.Cells.Item(168, 4).Value = 44452
$ws.Cells.Item(168, 5).Value = 10
$ws.Cells.Item(168, 6).Value = "Fruta"
$ws.Cells.Item(168, 7).Value = 100102
$ws.Cells.Item(168, 8).Value = "Cítricos"
$ws.Cells.Item(168, 9).Value = 100102003
$ws.Cells.Item(168, 10).Value = "Limón"
$ws.Cells.Item(168, 11).Value = "Sin especificar"
$ws.Cells.Item(168, 12).Value = "1a amarillo"
$ws.Cells.Item(168, 13).Value = 200
$ws.Cells.Item(168, 14).Value = 8500
$ws.Cells.Item(168, 15).Value = 8500
$ws.Cells.Item(168, 16).Value = 8500
$ws.Cells.Item(168, 17).Value = "`$/malla 16 kilos"
$ws.Cells.Item(168, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(168, 19).Value = 531
$ws.Cells.Item(168, 20).Value = 16

# New row 169: 2a amarillo
$ws.Cells.Item(169, 1).Value = 4
$ws.Cells.Item(169, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(169, 3).Value = "Los Lagos"
$ws.Cells.Item(169, 4).Value = 44452
$ws.Cells.Item(169, 5).Value = 10
$ws.Cells.Item(169, 6).Value = "Fruta"
$ws.Cells.Item(169, 7).Value = 100102
$ws.Cells.Item(169, 8).Value = "Cítricos"
$ws.Cells.Item(169, 9).Value = 100102003
$ws.Cells.Item(169, 10).Value = "Limón"
$ws.Cells.Item(169, 11).Value = "Sin especificar"
$ws.Cells.Item(169, 12).Value = "2a amarillo"
$ws.Cells.Item(169, 13).Value = 200
$ws.Cells.Item(169, 14).Value = 7500
$ws.Cells.Item(169, 15).Value = 7500
$ws.Cells.Item(169, 16).Value = 7500
$ws.Cells.Item(169, 17).Value = "`$/malla 16 kilos"
$ws.Cells.Item(169, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(169, 19).Value = 469
$ws.Cells.Item(169, 20).Value = 16
